$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# Insert a new worksheet column at J using EntireColumn.Insert
$ws.Range("J1").EntireColumn.Insert()

# Now resize table to include the new column S (since everything shifted right by 1)
$tbl.Resize($ws.Range("A1:S20"))

# Set header cell directly
$ws.Cells.Item(1, 10).Value = "archived"

# re-set all header cells to force sync (even to the same value, for shifted ones)
$ws.Cells.Item(1, 11).Value = "updating_each"
$ws.Cells.Item(1, 12).Value = "start_date"
$ws.Cells.Item(1, 13).Value = "end_date"
$ws.Cells.Item(1, 14).Value = "localisation"
$ws.Cells.Item(1, 15).Value = "delivery_format"
$ws.Cells.Item(1, 16).Value = "link"
$ws.Cells.Item(1, 17).Value = "data_path"
$ws.Cells.Item(1, 18).Value = "tag_ids"
$ws.Cells.Item(1, 19).Value = "doc_ids"

Write-Host "Table columns count:" $tbl.ListColumns.Count
for ($i = 1; $i -le $tbl.ListColumns.Count; $i++) {
    $col = $tbl.ListColumns.Item($i)
    Write-Host $i ":" $col.Name
}
